$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.939.14'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '3.805.79'
$ws.Range("E3").Value = '  -1.47%  '
$r = $ws.Range("D4")
$r.NumberFormat = '@'
$r.Value = '1.00'
$r.ClearFormats()
$ws.Range("E4").Value = '  -0.07%  '
$r = $ws.Range("D5")
$r.NumberFormat = '@'
$r.Value = '702.88'
$r.ClearFormats()
$ws.Range("E5").Value = '  +0.96%  '
$r = $ws.Range("D6")
$r.NumberFormat = '@'
$r.Value = '170.33'
$r.ClearFormats()
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("D7").Value = '3.805.61'
$ws.Range("E7").Value = '  -1.42%  '
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("E10").Value = '  -1.04%  '
$r = $ws.Range("D11")
$r.NumberFormat = '@'
$r.Value = '7.63'
$r.ClearFormats()
$ws.Range("E11").Value = '  +5.60%  '
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("E13").Value = '  -3.04%  '
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").Value = '4.447.85'
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").Value = '3.797.64'
$ws.Range("E16").Value = '  -1.90%  '
$ws.Range("D17").Value = '70.921.27'
$ws.Range("E17").Value = '  -0.26%  '
$r = $ws.Range("D18")
$r.NumberFormat = '@'
$r.Value = '17.36'
$r.ClearFormats()
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("E20").Value = '  -1.61%  '
$r = $ws.Range("D21")
$r.NumberFormat = '@'
$r.Value = '499.45'
$r.ClearFormats()
$ws.Range("E21").Value = '  +0.47%  '
$r = $ws.Range("D22")
$r.NumberFormat = '@'
$r.Value = '10.66'
$r.ClearFormats()
$ws.Range("E22").Value = '  -1.36%  '
$r = $ws.Range("D23")
$r.NumberFormat = '@'
$r.Value = '0.722'
$r.ClearFormats()
$ws.Range("E23").Value = '  +0.20%  '
$r = $ws.Range("D24")
$r.NumberFormat = '@'
$r.Value = '84.09'
$r.ClearFormats()
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("E25").Value = '  -5.20%  '
$ws.Range("D26").Value = '3.955.55'
$ws.Range("E26").Value = '  -1.23%  '
$r = $ws.Range("D27")
$r.NumberFormat = '@'
$r.Value = '12.04'
$r.ClearFormats()
$ws.Range("E27").Value = '  -1.41%  '
$r = $ws.Range("D28")
$r.NumberFormat = '@'
$r.Value = '10.28'
$r.ClearFormats()
$ws.Range("E28").Value = '  -3.70%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  -5.66%  '
$ws.Range("E31").Value = '  -4.20%  '
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("E33").Value = '  -3.32%  '
$r = $ws.Range("D34")
$r.NumberFormat = '@'
$r.Value = '29.00'
$r.ClearFormats()
$ws.Range("E34").Value = '  -2.06%  '
$ws.Range("E35").Value = '  -5.25%  '
$ws.Range("D36").Value = '3.771.55'
$ws.Range("E36").Value = '  -1.20%  '
$r = $ws.Range("D37")
$r.NumberFormat = '@'
$r.Value = '0.999'
$r.ClearFormats()
$ws.Range("E37").Value = '  -0.02%  '
$r = $ws.Range("D38")
$r.NumberFormat = '@'
$r.Value = '9.03'
$r.ClearFormats()
$ws.Range("E38").Value = '  -1.95%  '
$ws.Range("E39").Value = '  -3.34%  '
$r = $ws.Range("D40")
$r.NumberFormat = '@'
$r.Value = '2.36'
$r.ClearFormats()
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("E41").Value = '  -1.17%  '
$ws.Range("E42").Value = '  -3.02%  '
$ws.Range("E43").Value = '  -5.15%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  +0.07%  '
$r = $ws.Range("D46")
$r.NumberFormat = '@'
$r.Value = '166.78'
$r.ClearFormats()
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("E47").Value = '  +1.12%  '
$r = $ws.Range("D48")
$r.NumberFormat = '@'
$r.Value = '49.02'
$r.ClearFormats()
$ws.Range("E48").Value = '  -0.53%  '
$r = $ws.Range("D49")
$r.NumberFormat = '@'
$r.Value = '419.25'
$r.ClearFormats()
$ws.Range("E49").Value = '  +0.34%  '
$r = $ws.Range("D50")
$r.NumberFormat = '@'
$r.Value = '8.57'
$r.ClearFormats()
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("E51").Value = '  -2.75%  '
